$wb = $excel.ActiveWorkbook

$stage   = $wb.Worksheets.Item("STAGE")
$lmsprod = $wb.Worksheets.Item("LMSPROD")

# --- LMSPROD: new school / classroom / section + refreshed numeric ids ---
$lmsprod.Range("A2").Value = "FPK12School22538"
$lmsprod.Range("B2").Value = "FPK12Classroom3592"
$lmsprod.Range("C2").Value = "FPK12Section67301"

$lmsprod.Range("E3").NumberFormat = "@"
$lmsprod.Range("E3").Value = "57896"

$lmsprod.Range("E4").NumberFormat = "@"
$lmsprod.Range("E4").Value = "14780"

$lmsprod.Range("E5").NumberFormat = "@"
$lmsprod.Range("E5").Value = "53513"

# --- STAGE: refresh the tenant id used for this run, emphasised in bold ---
$stage.Range("D2").Value = "fpdistrict"
$stage.Range("D2").Font.Bold = $true
$stage.Range("D2").Font.Size = 8
$stage.Range("D2").Font.Name = "Arial"

# --- Selection / active-sheet state left by the last edit on STAGE ---
$stage.Activate()
$stage.Range("D2").Select()
